$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper functions
# ---------------------------------------------------------------------------

# Find the first occurrence of $searchText starting at character position
# $searchStart (defaults to 0) and return the matching Range (or $null).
function Find-TextRange($searchText, $searchStart) {
    $docEnd = $d.Content.End
    $r = $d.Range($searchStart, $docEnd)
    $found = $r.Find.Execute($searchText)
    if (-not $found) {
        Write-Host "NOT FOUND:" $searchText
        return $null
    }
    return $r
}

# Remove a leading "open tag" (e.g. "{g0}") and a trailing "close tag"
# (e.g. "{/g1}") that surround the text found via $searchText. Both tags are
# deleted using the Delete() method (never by assigning .Text on a range
# that starts exactly at a run boundary) so that the host run is preserved
# instead of being silently merged into its neighbour.
# Returns the character offset immediately following the (now stripped)
# inner text - handy for chaining subsequent searches.
function Strip-Tags($searchText, $openTag, $closeTag, $searchStart) {
    $r = Find-TextRange $searchText $searchStart
    if ($r -eq $null) { return $null }
    $start = $r.Start
    $end = $r.End

    $closeStart = $end - $closeTag.Length
    $rClose = $d.Range($closeStart, $end)
    if ($rClose.Text -ne $closeTag) {
        Write-Host "WARN close tag mismatch, expected" $closeTag "got" $rClose.Text
    }
    $rClose.Delete()

    $rOpen = $d.Range($start, $start + $openTag.Length)
    if ($rOpen.Text -ne $openTag) {
        Write-Host "WARN open tag mismatch, expected" $openTag "got" $rOpen.Text
    }
    $rOpen.Delete()

    return ($start + ($closeStart - ($start + $openTag.Length)))
}

# Replace just the numeric portion of a tag such as "{g16}" or "{/g17}" or
# "{x6}" with a new numeric string, leaving braces/slash/letter untouched.
# $fullTag is the complete tag text (e.g. "{/g17}"), $prefixLen is the
# number of characters preceding the digits (e.g. 3 for "{/g"), $digitsLen
# is how many digit characters follow, and $newDigits is the replacement.
function Renumber-Tag($searchText, $prefixLen, $digitsLen, $newDigits, $searchStart) {
    $r = Find-TextRange $searchText $searchStart
    if ($r -eq $null) { return $null }
    $start = $r.Start
    $rDigits = $d.Range($start + $prefixLen, $start + $prefixLen + $digitsLen)
    $rDigits.Text = $newDigits
    return $start
}

# ---------------------------------------------------------------------------
# 1. {g14},{/g15}  ->  ,
# ---------------------------------------------------------------------------
Strip-Tags "{g14},{/g15}" "{g14}" "{/g15}" 0

# ---------------------------------------------------------------------------
# 2. {g16} allhay apscay{/g17}  ->  {g14} allhay apscay{/g15}
# ---------------------------------------------------------------------------
Renumber-Tag "{g16}" 2 2 "14" 0
Renumber-Tag "{/g17}" 3 2 "15" 0

# ---------------------------------------------------------------------------
# 3. {g18}imestay ewnay omanray{/g19}  ->  {g16}imestay ewnay omanray{/g17}
# ---------------------------------------------------------------------------
Renumber-Tag "{g18}" 2 2 "16" 0
Renumber-Tag "{/g19}" 3 2 "17" 0

# ---------------------------------------------------------------------------
# 4. {g20}Arialhay, {/g21}  ->  {g18}Arialhay, {/g19}
# ---------------------------------------------------------------------------
Renumber-Tag "{g20}" 2 2 "18" 0
Renumber-Tag "{/g21}" 3 2 "19" 0

# ---------------------------------------------------------------------------
# 5. {g22}Arialhay 8 ptay{/g23}  ->  {g20}Arialhay 8 ptay{/g21}
# ---------------------------------------------------------------------------
Renumber-Tag "{g22}" 2 2 "20" 0
Renumber-Tag "{/g23}" 3 2 "21" 0

# ---------------------------------------------------------------------------
# 6. {g24}edray oregroundfay{/g25}  ->  {g22}edray oregroundfay{/g23}
# ---------------------------------------------------------------------------
Renumber-Tag "{g24}" 2 2 "22" 0
Renumber-Tag "{/g25}" 3 2 "23" 0

# ---------------------------------------------------------------------------
# 7. {g26},{/g27}  ->  ,
# ---------------------------------------------------------------------------
Strip-Tags "{g26},{/g27}" "{g26}" "{/g27}" 0

# ---------------------------------------------------------------------------
# 8. {g28} {/g29}  ->  {g24} {/g25}
# ---------------------------------------------------------------------------
Renumber-Tag "{g28}" 2 2 "24" 0
Renumber-Tag "{/g29}" 3 2 "25" 0

# ---------------------------------------------------------------------------
# 9. {g30}ueblay{/g31}  ->  {g26}ueblay{/g27}
# ---------------------------------------------------------------------------
Renumber-Tag "{g30}" 2 2 "26" 0
Renumber-Tag "{/g31}" 3 2 "27" 0

# ---------------------------------------------------------------------------
# 10. {g32},{/g33}  ->  ,
# ---------------------------------------------------------------------------
Strip-Tags "{g32},{/g33}" "{g32}" "{/g33}" 0

# ---------------------------------------------------------------------------
# 11. {g34} eengray, {/g35}  ->  {g28} eengray, {/g29}
# ---------------------------------------------------------------------------
Renumber-Tag "{g34}" 2 2 "28" 0
Renumber-Tag "{/g35}" 3 2 "29" 0

# ---------------------------------------------------------------------------
# 12. {g36}ellowyay ighlighthay{/g37}  ->  {g30}ellowyay ighlighthay{/g31}
# ---------------------------------------------------------------------------
Renumber-Tag "{g36}" 2 2 "30" 0
Renumber-Tag "{/g37}" 3 2 "31" 0

# ---------------------------------------------------------------------------
# 13. {x38}  ->  {x32}   (inside first hyperlink)
# ---------------------------------------------------------------------------
Renumber-Tag "{x38}" 2 2 "32" 0

# ---------------------------------------------------------------------------
# 14. {g39}erlinkhypay{/g40}  ->  {g33}erlinkhypay{/g34}
# ---------------------------------------------------------------------------
Renumber-Tag "{g39}" 2 2 "33" 0
Renumber-Tag "{/g40}" 3 2 "34" 0

# ---------------------------------------------------------------------------
# 15. {x41}, ahay   ->  {x35}, ahay
# ---------------------------------------------------------------------------
Renumber-Tag "{x41}" 2 2 "35" 0

# ---------------------------------------------------------------------------
# 16. {x42}  ->  {x36}   (inside second hyperlink)
# ---------------------------------------------------------------------------
Renumber-Tag "{x42}" 2 2 "36" 0

# ---------------------------------------------------------------------------
# 17. {g43}ookmarkbay umpjay otay ethay orderedhay istlay{/g44}
#      ->  {g37}ookmarkbay umpjay otay ethay orderedhay istlay{/g38}
# ---------------------------------------------------------------------------
Renumber-Tag "{g43}" 2 2 "37" 0
Renumber-Tag "{/g44}" 3 2 "38" 0

# ---------------------------------------------------------------------------
# 18. {x45}, andhay ahay ootnotefay  ->  {x39}, andhay ahay ootnotefay
# ---------------------------------------------------------------------------
Renumber-Tag "{x45}" 2 2 "39" 0

# ---------------------------------------------------------------------------
# 19. {x46} + {g47}.{/g48}  (two runs)  ->  {x40}.  (single run)
# ---------------------------------------------------------------------------
Renumber-Tag "{x46}" 2 2 "40" 0
Strip-Tags "{g47}.{/g48}" "{g47}" "{/g48}" 0

# ---------------------------------------------------------------------------
# 20. {x3} + {g4}ommentscay{/g5}  (two runs)  ->  {x3}ommentscay  (single run)
#     Note: {x3} keeps its number.
# ---------------------------------------------------------------------------
Strip-Tags "{g4}ommentscay{/g5}" "{g4}" "{/g5}" 0

# ---------------------------------------------------------------------------
# 21. {x6}  ->  {x4}
# ---------------------------------------------------------------------------
Renumber-Tag "{x6}" 2 1 "4" 0

# ---------------------------------------------------------------------------
# 22. {x7} + {g8}.{/g9}  (two runs)  ->  {x5}.  (single run)
# ---------------------------------------------------------------------------
Renumber-Tag "{x7}" 2 1 "5" 0
Strip-Tags "{g8}.{/g9}" "{g8}" "{/g9}" 0

# ---------------------------------------------------------------------------
# 23-42. Plain {g0}...{/g1} tag pairs that are simply stripped.
# ---------------------------------------------------------------------------
Strip-Tags "{g0}isthay aragraphpay ashay ahay ueblay outlinehay.{/g1}" "{g0}" "{/g1}" 0
Strip-Tags "{g0}isthay ishay anhay orderedhay istlay:{/g1}" "{g0}" "{/g1}" 0
Strip-Tags "{g0}Onehay{/g1}" "{g0}" "{/g1}" 0
Strip-Tags "{g0}otway{/g1}" "{g0}" "{/g1}" 0
Strip-Tags "{g0}eethray{/g1}" "{g0}" "{/g1}" 0
Strip-Tags "{g0}isthay ishay anhay unorderedhay istlay:{/g1}" "{g0}" "{/g1}" 0
Strip-Tags "{g0}Applehay{/g1}" "{g0}" "{/g1}" 0
Strip-Tags "{g0}acintoshmay{/g1}" "{g0}" "{/g1}" 0

# "onagoldjay" strip + the {x2} that immediately follows becomes {x0}
$afterOnagoldjay = Strip-Tags "{g0}onagoldjay{/g1}" "{g0}" "{/g1}" 0
Renumber-Tag "{x2}" 2 1 "0" $afterOnagoldjay

Strip-Tags "{g0}ananabay{/g1}" "{g0}" "{/g1}" 0
Strip-Tags "{g0}Orangehay{/g1}" "{g0}" "{/g1}" 0
Strip-Tags "{g0}Ahay abletay ollowsfay:{/g1}" "{g0}" "{/g1}" 0
Strip-Tags "{g0}olumncay 1 owray 1{/g1}" "{g0}" "{/g1}" 0
Strip-Tags "{g0}cay2ray1{/g1}" "{g0}" "{/g1}" 0
Strip-Tags "{g0}cay1ray2{/g1}" "{g0}" "{/g1}" 0
Strip-Tags "{g0}cay2ray2{/g1}" "{g0}" "{/g1}" 0
Strip-Tags "{g0}ollowingfay ishay ahay anualmay agepay eakbray:{/g1}" "{g0}" "{/g1}" 0
Strip-Tags "{g0}2008 anjay 23: irstfay ersionvay{/g1}" "{g0}" "{/g1}" 0
Strip-Tags "{g0}2008 ebfay 06: vay1.1: ixedfay irstfay erlinkhypay; angedchay Excelhay omfray egularray astepay ashay abletay otay anhay embeddedhay eadsheetspray{/g1}" "{g0}" "{/g1}" 0

Write-Host "Done."
